$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.199.62"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "1.658.52"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.56"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5216"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2671"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.15"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.442"
$ws.Range("D13").Value = "1.638.50"
$ws.Range("E13").Value = "  -2.48%  "
$ws.Range("D14").Value = "1.884.30"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5490"
$ws.Range("E15").Value = "  -2.51%  "
$ws.Range("D16").Value = "0.0₅8265"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "26.248.85"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.682"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.41"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.118"
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.08"
$ws.Range("E25").Value = "  -3.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1246"
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.257"
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.19"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06022"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.285"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.344"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.652"
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9842"
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.411"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.771"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5938"
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.968"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8641"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "1.041.68"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.77"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "1.798.97"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.117"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05179"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.471"
$ws.Range("E51").Value = "  +3.88%  "
